$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting existing Age (E) / Sex (F) columns
# to F / G respectively.
$ws.Columns("E").Insert()

# Set the new header for column E
$ws.Cells.Item(1, 5).Value = "StimOrder"

# Fill in the StimOrder values for each subject row (2-49)
$stimOrder = @(321, 312, 132, 321, 321, 123, 231, 231, 231, 213, 321, 213, 213, 132, 123, 123, 321, 213, 123, 123, 312, 213, 123, 321, 231, 312, 321, 213, 132, 312, 132, 123, 213, 132, 312, 123, 312, 213, 132, 231, 231, 231, 312, 312, 132, 123, 321, 231)

for ($i = 0; $i -lt $stimOrder.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $stimOrder[$i]
}
